$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.272.53'
$ws.Range("E2").Value = '  +1.88%  '
$ws.Range("D3").Value = '2.085.99'
$ws.Range("E3").Value = '  -0.51%  '
$ws.Range("E4").Value = '  -0.57%  '
$ws.Range("D5").Value = "'342.98"
$ws.Range("E5").Value = '  -0.54%  '
$ws.Range("E6").Value = '  -0.46%  '
$ws.Range("D7").Value = "'0.5229"
$ws.Range("E7").Value = '  +1.58%  '
$ws.Range("D8").Value = "'0.4407"
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").Value = "'54.50"
$ws.Range("E9").Value = '  +3.22%  '
$ws.Range("E10").Value = '  +0.36%  '
$ws.Range("E11").Value = '  -0.44%  '
$ws.Range("E12").Value = '  -0.56%  '
$ws.Range("D13").Value = "'8.558"
$ws.Range("E13").Value = '  +3.04%  '
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").Value = "'6.885"
$ws.Range("E14").Value = '  +1.92%  '
$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").Value = '2.023.15'
$ws.Range("E15").Value = '  -4.37%  '
$ws.Range("D16").Value = "'101.16"
$ws.Range("E16").Value = '  +1.47%  '
$ws.Range("D17").Value = "'0.00001155"
$ws.Range("E17").Value = '  +0.31%  '
$ws.Range("D18").Value = "'1.004"
$ws.Range("E18").Value = '  -0.40%  '
$ws.Range("D19").Value = "'21.03"
$ws.Range("E19").Value = '  +0.73%  '
$ws.Range("E20").Value = '  +0.47%  '
$ws.Range("E21").Value = '  +1.90%  '
$ws.Range("D22").Value = "'1.003"
$ws.Range("E22").Value = '  -0.36%  '
$ws.Range("D23").Value = '30.320.71'
$ws.Range("E23").Value = '  +1.86%  '
$ws.Range("E24").Value = '  -0.86%  '
$ws.Range("D25").Value = "'2.306"
$ws.Range("E25").Value = '  -0.28%  '
$ws.Range("D26").Value = "'21.74"
$ws.Range("E26").Value = '  -0.70%  '
$ws.Range("D27").Value = "'162.29"
$ws.Range("E27").Value = '  +0.15%  '
$ws.Range("D28").Value = "'2.501"
$ws.Range("D29").Value = "'132.86"
$ws.Range("E29").Value = '  -0.11%  '
$ws.Range("D30").Value = "'1.128"
$ws.Range("E30").Value = '  -0.15%  '
$ws.Range("E31").Value = '  +0.64%  '
$ws.Range("E32").Value = '  -0.77%  '
$ws.Range("D33").Value = "'6.209"
$ws.Range("E33").Value = '  +0.57%  '
$ws.Range("D34").Value = "'6.675"
$ws.Range("E34").Value = '  +9.04%  '
$ws.Range("D35").Value = "'3.855"
$ws.Range("E35").Value = '  -2.28%  '
$ws.Range("D36").Value = "'10.22"
$ws.Range("E36").Value = '  -1.78%  '
$ws.Range("D37").Value = "'0.02618"
$ws.Range("E37").Value = '  +1.86%  '
$ws.Range("D38").Value = "'0.06813"
$ws.Range("E38").Value = '  +1.34%  '
$ws.Range("B39").Value = 'TrustWalletToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D39").Value = "'1.347"
$ws.Range("E39").Value = '  +3.83%  '
$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D40").Value = "'0.6951"
$ws.Range("E40").Value = '  +1.38%  '
$ws.Range("D41").Value = "'12.46"
$ws.Range("E41").Value = '  +0.01%  '
$ws.Range("D42").Value = "'0.2202"
$ws.Range("E42").Value = '  -1.12%  '
$ws.Range("D43").Value = "'0.6788"
$ws.Range("E43").Value = '  +2.30%  '
$ws.Range("D44").Value = "'14.30"
$ws.Range("E44").Value = '  +0.02%  '
$ws.Range("B45").Value = 'Frax'
$ws.Range("C45").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D45").Value = "'1.003"
$ws.Range("E45").Value = '  -0.39%  '
$ws.Range("B46").Value = 'NEARProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D46").Value = "'2.317"
$ws.Range("E46").Value = '  -0.18%  '
$ws.Range("E47").Value = '  +17.80%  '
$ws.Range("D48").Value = "'3.640"
$ws.Range("E48").Value = '  +0.61%  '
$ws.Range("D49").Value = "'0.00000000348"
$ws.Range("E49").Value = '  +0.41%  '
$ws.Range("D50").Value = "'1.203"
$ws.Range("E50").Value = '  +7.52%  '
$ws.Range("D51").Value = "'1.213"
$ws.Range("E51").Value = '  -0.72%  '
